$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value for column D (Price), kept as TEXT
# (the source data stores these as inline strings, not numbers).
$priceUpdates = @{
    "D2"  = "244.22"
    "D3"  = "23.95"
    "D4"  = "5.259"
    "D6"  = "6.461"
    "D7"  = "3.266"
    "D9"  = "0.8858"
    "D10" = "0.1381"
    "D11" = "0.07132"
    "D12" = "0.03079"
    "D13" = "0.03057"
    "D14" = "0.09326"
    "D15" = "3.813"
    "D16" = "0.001541"
    "D17" = "0.04712"
    "D18" = "0.0006010"
    "D19" = "0.006180"
    "D22" = "0.00008700"
    "D23" = "3.545"
    "D24" = "2.170"
    "D26" = "0.1312"
    "D40" = "0.03840"
    "D41" = "0.006270"
    "D42" = "0.1053"
    "D43" = "0.002542"
    "D44" = "0.007282"
    "D45" = "0.00005330"
    "D47" = "0.5500"
    "D48" = "0.003777"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Column E text tweaks (the "Worstin24h" tag moved from row 18 to row 47)
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
